$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The model architecture value in B2 was renamed from "sep_cls" to "cls_sep"
$ws.Range("B2").Value = "cls_sep"

# Update the active selection to a single cell B3 (was F2:F7 with active cell F2)
$ws.Range("B3").Select()
